$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $r = $ws.Range($addr)
    $origStyle = $r.Style
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = $origStyle
}

Set-TextValue "D2" '67.807.09'
Set-TextValue "E2" '  -0.10%  '

Set-TextValue "D3" '3.772.36'
Set-TextValue "E3" '  -1.61%  '

Set-TextValue "E4" '  +0.30%  '

Set-TextValue "D5" '604.70'
Set-TextValue "E5" '  +0.28%  '

Set-TextValue "D6" '162.52'
Set-TextValue "E6" '  -2.57%  '

Set-TextValue "D7" '3.769.56'
Set-TextValue "E7" '  -1.61%  '

Set-TextValue "E8" '  +0.14%  '

Set-TextValue "E10" '  -2.50%  '

Set-TextValue "B11" 'Toncoin'
Set-TextValue "C11" 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextValue "D11" '6.82'
Set-TextValue "E11" '  +8.28%  '

Set-TextValue "B12" 'Cardano'
Set-TextValue "C12" 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
Set-TextValue "D12" '0.445'
Set-TextValue "E12" '  -1.52%  '

Set-TextValue "D13" '0.0000244'
Set-TextValue "E13" '  -3.67%  '

Set-TextValue "D14" '34.79'
Set-TextValue "E14" '  -3.19%  '

Set-TextValue "D15" '4.404.01'
Set-TextValue "E15" '  -1.47%  '

Set-TextValue "D16" '3.779.26'
Set-TextValue "E16" '  -1.71%  '

Set-TextValue "D17" '67.809.57'
Set-TextValue "E17" '  -0.16%  '

Set-TextValue "D18" '17.99'
Set-TextValue "E18" '  -2.32%  '

Set-TextValue "E19" '  +1.93%  '

Set-TextValue "D20" '6.97'
Set-TextValue "E20" '  -1.53%  '

Set-TextValue "D21" '456.31'
Set-TextValue "E21" '  -1.78%  '

Set-TextValue "D22" '9.39'
Set-TextValue "E22" '  -5.50%  '

Set-TextValue "D23" '0.687'
Set-TextValue "E23" '  -2.08%  '

Set-TextValue "D24" '0.0000147'
Set-TextValue "E24" '  -1.23%  '

Set-TextValue "D25" '82.99'
Set-TextValue "E25" '  -0.58%  '

Set-TextValue "D26" '11.84'
Set-TextValue "E26" '  -2.06%  '

Set-TextValue "B27" 'Dai'
Set-TextValue "C27" 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextValue "D27" '1.00'
Set-TextValue "E27" '  -0.01%  '

Set-TextValue "B28" 'Fetch.AI'
Set-TextValue "C28" 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
Set-TextValue "D28" '2.06'
Set-TextValue "E28" '  -3.04%  '

Set-TextValue "D29" '9.88'
Set-TextValue "E29" '  -1.97%  '

Set-TextValue "D30" '3.914.62'
Set-TextValue "E30" '  -1.66%  '

Set-TextValue "D31" '2.59'
Set-TextValue "E31" '  -7.05%  '

Set-TextValue "D32" '7.14'
Set-TextValue "E32" '  -3.63%  '

Set-TextValue "D33" '2.15'
Set-TextValue "E33" '  -3.24%  '

Set-TextValue "D34" '28.78'
Set-TextValue "E34" '  -2.90%  '

Set-TextValue "D35" '0.998'
Set-TextValue "E35" '  -0.21%  '

Set-TextValue "D36" '8.87'
Set-TextValue "E36" '  -2.77%  '

Set-TextValue "D37" '0.0990'
Set-TextValue "E37" '  -1.27%  '

Set-TextValue "E38" '  +7.14%  '

Set-TextValue "D39" '5.80'
Set-TextValue "E39" '  -0.37%  '

Set-TextValue "B40" 'dogwifhat'
Set-TextValue "C40" 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
Set-TextValue "D40" '3.21'
Set-TextValue "E40" '  -3.25%  '

Set-TextValue "B41" 'Mantle'
Set-TextValue "C41" 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
Set-TextValue "D41" '0.974'
Set-TextValue "E41" '  -2.71%  '

Set-TextValue "D42" '0.999'
Set-TextValue "E42" '  -0.05%  '

Set-TextValue "B44" 'OKB'
Set-TextValue "C44" 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue "D44" '47.04'
Set-TextValue "E44" '  -2.06%  '

Set-TextValue "B45" 'Monero'
Set-TextValue "C45" 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue "D45" '152.59'
Set-TextValue "E45" '  +2.00%  '

Set-TextValue "B46" 'Arweave'
Set-TextValue "C46" 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
Set-TextValue "D46" '43.01'
Set-TextValue "E46" '  -3.58%  '

Set-TextValue "B47" 'ONDO'
Set-TextValue "C47" 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
Set-TextValue "D47" '1.39'
Set-TextValue "E47" '  -3.06%  '

Set-TextValue "B48" 'TheGraph'
Set-TextValue "C48" 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
Set-TextValue "D48" '0.291'
Set-TextValue "E48" '  -3.67%  '

Set-TextValue "D49" '8.28'
Set-TextValue "E49" '  -0.71%  '

Set-TextValue "D50" '1.82'
Set-TextValue "E50" '  -1.29%  '

Set-TextValue "D51" '26.23'
Set-TextValue "E51" '  -9.14%  '
